$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 99, shifting existing rows 99-184 down to 100-185
$ws.Rows.Item(99).Insert()

# Populate the newly inserted row 99 with the weekly data point
$ws.Range("A99").Value = 5
$ws.Range("B99").Value = "Macroferia Regional de Talca"
$ws.Range("C99").Value = "Maule"
$ws.Range("D99").Value = 44447
$ws.Range("E99").Value = 7
$ws.Range("F99").Value = 100112023
$ws.Range("G99").Value = "Brócoli"
$ws.Range("H99").Value = "Sin especificar"
$ws.Range("I99").Value = "Primera"
$ws.Range("J99").Value = 3000
$ws.Range("K99").Value = 600
$ws.Range("L99").Value = 600
$ws.Range("M99").Value = 600
$ws.Range("N99").Value = "$/unidad"
$ws.Range("O99").Value = "Región del Maule"
$ws.Range("P99").Value = 600
$ws.Range("Q99").Value = 1
$ws.Range("R99").Value = "Hortaliza"
